$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert new rows (top to bottom, using positions valid at each step) ---
$ws.Rows.Item(114).Insert()
$ws.Range("A122:A130").EntireRow.Insert()
$ws.Range("A136:A138").EntireRow.Insert()

# --- Step 2: Set cell values for newly inserted content ---
$ws.Range("C113").Value = "Forward Rate Models: Volatility Models"
$ws.Range("C122").Value = "Forward Rate Models: Volatility Models"
$ws.Range("C124").Value = "ShortRateVolatilityModelInterface"
$ws.Range("D124").Value = "ShortRateVolatilityModelInterface"
$ws.Range("E124").Value = "interface"
$ws.Range("C125").Value = "AbstractShortRateVolatilityModel"
$ws.Range("D125").Value = "AbstractShortRateVolatilityModel"
$ws.Range("E125").Value = "partial implementation"
$ws.Range("C127").Value = "AbstractShortRateVolatilityModelParametric"
$ws.Range("D127").Value = "ShortRateVolatilityModelParametric (extract interface)"
$ws.Range("E127").Value = "interface"
$ws.Range("C128").Value = "AbstractShortRateVolatilityModelParametric"
$ws.Range("D128").Value = "AbstractShortRateVolatilityModelParametric"
$ws.Range("E128").Value = "partial implementation"
$ws.Range("C129").Value = "ShortRateVolatilityModelCalibrateable"
$ws.Range("D129").Value = "ShortRateVolatilityModelCalibrateable"
$ws.Range("E129").Value = "interface"
$ws.Range("C136").Value = "double[] AbstractShortRateVolatilityModelParametric.getParameter()"
$ws.Range("D136").Value = "double[] AbstractShortRateVolatilityModelParametric.getParameterAsDouble()"
$ws.Range("E136").Value = "method"
$ws.Range("D137").Value = "RandomVariable[] AbstractShortRateVolatilityModelParametric.getParameter()"
$ws.Range("E137").Value = "method"

# --- Step 3: Apply styles to new cells by copying format from reference cells ---
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B113").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C113").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B122").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C122").PasteSpecial(-4122) | Out-Null
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A124").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E124").PasteSpecial(-4122) | Out-Null
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A125").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A127").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E127").PasteSpecial(-4122) | Out-Null
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A128").PasteSpecial(-4122) | Out-Null
$ws.Range("A49").Copy() | Out-Null
$ws.Range("A129").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E129").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A134").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A136").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A137").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 4: Update sheet view / selection ---
$ws.Range("B116").Select() | Out-Null
